# Added Configurable zero_before_threshold parameter to enable setting dims
# before noise_threshold or First Rise Point to 0.
#
# This updates the First_Noticeable_Increase_Index (C), the
# First_Noticeable_Increase_Cumulative_Value (E) and the Pulse_Width (G)
# columns on each of the Step3_DataPts_* sheets to reflect the new
# zero_before_threshold behaviour.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Step3_DataPts_0.5", "Step3_DataPts_0.7", "Step3_DataPts_0.8", "Step3_DataPts_0.9")

# Values common to every Step3_DataPts_* sheet (rows 2, 3 and 4 share the
# same First_Noticeable_Increase_Index / Cumulative_Value across thresholds).
$commonRowValues = @{
    2 = @{ C = 88; E = 0.002304090588778202 }
    3 = @{ C = 88; E = 0.00477131187798357 }
    4 = @{ C = 88; E = 0.007013136003042588 }
    6 = @{ C = 87; E = 0.004000218324993229 }
}

# Pulse_Width (column G) is unique per sheet/row combination.
$pulseWidthBySheet = @{
    "Step3_DataPts_0.5" = @{ 2 = 18; 3 = 18; 4 = 18; 6 = 19 }
    "Step3_DataPts_0.7" = @{ 2 = 29; 3 = 27; 4 = 24; 6 = 23 }
    "Step3_DataPts_0.8" = @{ 2 = 66; 3 = 62; 4 = 47; 6 = 30 }
    "Step3_DataPts_0.9" = @{ 2 = 80; 3 = 79; 4 = 68; 6 = 68 }
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $pulseWidths = $pulseWidthBySheet[$sheetName]

    foreach ($row in $commonRowValues.Keys) {
        $vals = $commonRowValues[$row]

        $ws.Range("C$row").Value = $vals.C
        $ws.Range("E$row").Value = $vals.E
        $ws.Range("G$row").Value = $pulseWidths[$row]
    }
}
